{"js": "// Target paragraph texts after the edit (28 paragraphs total).\n// A \"\\u000b\" (vertical tab) represents a line break (<w:br/>) within a run,\n// matching how Word/Office.js exposes in-paragraph breaks through Range.text.\nconst targetTexts = [\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 24.01.25\\u000bBack to Basics: Revisiting REINFORCE Style Optimization for Learning from Human Feedback in LLMs\",\n  \"\u05ea\u05de\u05e6\u05d9\u05ea \u05d4\u05de\u05d0\u05de\u05e8:\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05df \u05de\u05d7\u05d3\u05e9 \u05d0\u05ea \u05d4\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05dc\u05de\u05d9\u05d3\u05d4 \u05de\u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd \u05de\u05e4\u05d9\u05d3\u05d1\u05e7 \u05d0\u05e0\u05d5\u05e9\u05d9 (RLHF) \u05d1\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05ea LLMs. \u05d4\u05d5\u05d0 \u05de\u05d0\u05ea\u05d2\u05e8 \u05d0\u05ea \u05d4\u05d3\u05d5\u05de\u05d9\u05e0\u05e0\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc PPO (Proximal Policy Optimization) \u05db\u05e9\u05d9\u05d8\u05ea \u05dc\u05de\u05d9\u05d3\u05ea \u05d4\u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd \u05d4\u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea \u05d1\u05d4\u05e7\u05e9\u05e8 \u05d6\u05d4, \u05ea\u05d5\u05da \u05d4\u05d3\u05d2\u05e9\u05ea \u05d7\u05d5\u05e1\u05e8 \u05d4\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d5\u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d4\u05de\u05d9\u05d5\u05ea\u05e8\u05ea \u05e9\u05dc\u05d5. \u05d1\u05de\u05e7\u05d5\u05dd \u05d6\u05d0\u05ea, \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05dc\u05d7\u05d6\u05d5\u05e8 \u05dc\u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e1\u05d2\u05e0\u05d5\u05df REINFORCE, \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05ea Vanilla Policy Gradient (REINFORCE) \u05d5\u05d4\u05e8\u05d7\u05d1\u05ea\u05d5 \u05d4\u05e8\u05d1-\u05d3\u05d2\u05d9\u05de\u05ea\u05d9\u05ea, REINFORCE Leave-One-Out (RLOO). \u05e9\u05d9\u05d8\u05d5\u05ea \u05d0\u05dc\u05d5 \u05de\u05d5\u05db\u05d9\u05d7\u05d5\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de-PPO \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05e2\u05dc\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea, \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05de\u05e1\u05e4\u05e8 \u05de\u05e2\u05e8\u05db\u05d9 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05d5\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea LLM. \u05d4\u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05e9\u05d9\u05dd \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05e9\u05d9\u05d2 \u05d4\u05ea\u05d0\u05de\u05d4 \u05e9\u05dc LLMs \u05dc\u05d4\u05e2\u05d3\u05e4\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea \u05e2\u05dd \u05d0\u05e1\u05d8\u05e8\u05d8\u05d2\u05d9\u05d5\u05ea \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d4\u05de\u05d5\u05ea\u05d0\u05de\u05d5\u05ea \u05dc\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05d5\u05ea \u05e9\u05dc RLHF.\",\n  \"\u05d4\u05e8\u05d7\u05d1\u05d4 \u05e2\u05dc \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea:\",\n  \"1. \u05e4\u05d9\u05e9\u05d5\u05d8 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9:\",\n  \"\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05e8\u05d1\u05d9\u05dd \u05de\u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc PPO (\u05dc\u05de\u05e9\u05dc, \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2, \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05e8\u05da (value), \u05d5\u05de\u05d9\u05d3\u05d5\u05dc  \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd) \u05d0\u05d9\u05e0\u05dd \u05d4\u05db\u05e8\u05d7\u05d9\u05d9\u05dd \u05dc-RLHF, \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05d0\u05ea\u05d7\u05d5\u05dc \u05d8\u05d5\u05d1 \u05e9\u05dc LLM (\u05dc\u05d0\u05d7\u05e8 SFT \u05dc\u05de\u05e9\u05dc). \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d9\u05d3\u05d5\u05dc \u05e1\u05d3\u05e8\u05d5\u05ea \u05e9\u05dc\u05de\u05d5\u05e5 \u05db\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d5\u05d3\u05d3\u05d5\u05ea, REINFORCE \u05e0\u05de\u05e0\u05e2 \u05de\u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05e8\u05da-\u05de\u05e6\u05d1(V \u05d5-Q) \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05df, \u05d5\u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05d1\u05e2\u05d9\u05d4 \u05dc\u05d3\u05d5\u05de\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05dc\u05d1\u05e0\u05d3\u05d9\u05d8 \u05d4\u05e7\u05e9\u05e8\u05d9.\",\n  \"2. \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05de\u05e2\u05e9\u05d9\u05ea:\",\n  \"\u05e9\u05d9\u05d8\u05ea RLOO \u05de\u05e9\u05ea\u05de\u05e9\u05ea \u05d1\u05db\u05dc \u05d4\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05e9\u05e0\u05d5\u05e6\u05e8\u05d5 \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05d1\u05e1\u05d9\u05e1 \u05d4\u05e9\u05d5\u05d5\u05d0\u05d4, \u05de\u05e9\u05d9\u05d2 \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05d2\u05d1\u05d5\u05d4\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de-RAFT, \u05e9\u05dc\u05d5\u05e7\u05d7 \u05d0\u05ea \u05e8\u05e7 \u05d0\u05ea \u05d4\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e6\u05d9\u05d5\u05e0\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd (\u05e1\u05d5\u05d2 \u05e9\u05dc rejection sampling). \u05d6\u05d4 \u05de\u05d5\u05d1\u05d9\u05dc \u05dc\u05d7\u05d9\u05e1\u05db\u05d5\u05df \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9 \u05d1\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05dd \u05d5\u05e0\u05d9\u05e6\u05d5\u05dc \u05d8\u05d5\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05d4\u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05d4\u05d6\u05de\u05d9\u05e0\u05d9\u05dd. \u05d4\u05d2\u05d9\u05e9\u05d4 \u05de\u05e4\u05e9\u05d8\u05ea \u05d0\u05ea \u05ea\u05d4\u05dc\u05d9\u05db\u05d9 \u05d4-RLHF \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05e4\u05d7\u05ea\u05ea \u05d4\u05ea\u05dc\u05d5\u05ea \u05d1\u05d4\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e8\u05d2\u05d9\u05e9\u05d9\u05dd \u05db\u05de\u05d5 \u05d9\u05d7\u05e1\u05d9 \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05d5\u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d1\u05e9\u05e2\u05e8\u05d5\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d9\u05ea\u05e8\u05d5\u05df (\u05db\u05de\u05d5 \u05d1-GAE).\",\n  \"4.\u05e8\u05d5\u05d1\u05e1\u05d8\u05d9\u05d5\u05ea:\\u000b\u05e9\u05d9\u05d8\u05ea RLOO  \u05de\u05d3\u05d2\u05d9\u05de\u05d4 \u05e8\u05d5\u05d1\u05e1\u05d8\u05d9\u05d5\u05ea \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05e8\u05d5\u05e2\u05e9\u05d9\u05dd \u05d5\u05e2\u05d5\u05e0\u05e9\u05d9 KL \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8, \u05e2\u05d5\u05dc\u05d4 \u05e2\u05dc \u05e9\u05d9\u05d8\u05d5\u05ea \u05db\u05de\u05d5 RAFT \u05e9\u05e8\u05d2\u05d9\u05e9\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05dc\u05d3\u05d9\u05d5\u05e7\u05dd.\",\n  \"\u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea:\",\n  \"1. \u05d0\u05d9\u05d6\u05d5\u05df \u05e9\u05d5\u05e0\u05d5\u05ea-\u05d4\u05d8\u05d9\u05d4(bias-variance tradeoff) \u05d5\u05e9\u05e2\u05e8\u05d5\u05da \u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05dc\u05dc\u05d0 \u05d4\u05d8\u05d9\u05d4:\",\n  \"\u05e9\u05d9\u05d8\u05ea PPO \u05de\u05e1\u05ea\u05de\u05db\u05ea \u05e2\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05e8\u05da-\u05de\u05e6\u05d1 \u05d5\u05e9\u05e2\u05e8\u05d5\u05da \u05d9\u05ea\u05e8\u05d5\u05df \u05de\u05d5\u05db\u05dc\u05dc (Generalized Advantage Estimation) \u05dc\u05d4\u05e4\u05d7\u05ea\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05e9\u05e2\u05e8\u05d5\u05da \u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05d1\u05de\u05d7\u05d9\u05e8 \u05e9\u05dc \u05d4\u05e2\u05dc\u05d0\u05ea \u05d4\u05d8\u05d9\u05d4. \u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05e9\u05d1-RLHF \u05e2\u05d1\u05d5\u05e8 LLMs \u05d3\u05d9 \u05de\u05d0\u05d5\u05de\u05df (warm start) \u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05e4\u05d7\u05ea\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05dc\u05e4\u05d7\u05d5\u05ea \u05e7\u05e8\u05d9\u05d8\u05d9\u05ea. \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e9\u05d9\u05d8\u05d5\u05ea \u05dc\u05dc\u05d0 \u05d4\u05d8\u05d9\u05d4 \u05db\u05de\u05d5 REINFORCE \u05dc\u05ea\u05e4\u05e7\u05d3 \u05d4\u05d9\u05d8\u05d1 \u05d1\u05dc\u05d9 \u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05d4\u05d8\u05d9\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea. \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea, \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d3\u05d2\u05d9\u05dd \u05e9-REINFORCE \u05de\u05e9\u05d9\u05d2 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05d8\u05d5\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de-PPO, \u05d0\u05e4\u05d9\u05dc\u05d5 \u05ea\u05d7\u05ea \u05ea\u05e0\u05d0\u05d9\u05dd \u05e9\u05dc \u05e9\u05d5\u05e0\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4\u05d4 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea.\",\n  \"2. \u05de\u05d9\u05d3\u05d5\u05dc \u05de\u05e1\u05dc\u05d5\u05dc \u05de\u05dc\u05d0(\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e9\u05dc\u05de\u05d4) \u05dc\u05e2\u05d5\u05de\u05ea \u05de\u05d9\u05d3\u05d5\u05dc \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd:\",\n  \"\u05e9\u05d9\u05d8\u05ea PPO \u05de\u05de\u05d3\u05dc\u05ea \u05db\u05dc \u05d8\u05d5\u05e7\u05df \u05db\u05e4\u05e2\u05d5\u05dc\u05d4, \u05d9\u05d5\u05e6\u05e8 \u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d7\u05dc\u05d8\u05d4 \u05de\u05e8\u05e7\u05d5\u05d1\u05d9 (MDP) \u05d1\u05d5 \u05e8\u05e6\u05e4\u05d9 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d7\u05dc\u05e7\u05d9\u05d9\u05dd \u05d4\u05dd \u05de\u05e6\u05d1\u05d9\u05dd. \u05e2\u05dd \u05d6\u05d0\u05ea, RLHF \u05de\u05d9\u05d9\u05d7\u05e1 \u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05e8\u05e7 \u05dc\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05e9\u05dc\u05de\u05d5\u05ea, \u05de\u05d4 \u05e9\u05d4\u05d5\u05e4\u05da \u05de\u05e6\u05d1\u05d9 \u05d1\u05d9\u05e0\u05d9\u05d9\u05dd \u05dc\u05dc\u05d0 \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d9\u05dd. \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d9\u05d3\u05d5\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05db\u05e4\u05e2\u05d5\u05dc\u05d4 \u05d9\u05d7\u05d9\u05d3\u05d4, REINFORCE \u05de\u05e4\u05e9\u05d8 \u05d0\u05ea \u05d4\u05d1\u05e2\u05d9\u05d4 \u05dc\u05de\u05d1\u05e0\u05d4 Contextual Bandit \u05d4\u05de\u05ea\u05d9\u05d9\u05e9\u05e8 \u05d9\u05e9\u05d9\u05e8\u05d5\u05ea \u05e2\u05dd \u05de\u05d1\u05e0\u05d4 \u05d4\u05ea\u05d2\u05de\u05d5\u05dc. \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea \u05de\u05d0\u05e9\u05e8\u05d5\u05ea \u05db\u05d9 \u05d2\u05d9\u05e9\u05d4 \u05d6\u05d5 \u05e2\u05d5\u05dc\u05d4 \u05e2\u05dc \u05de\u05d9\u05d3\u05d5\u05dc \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05df \u05d1\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d5\u05d4\u05df \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd.\",\n  \"3. \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05d5\u05d9\u05e6\u05d9\u05d1\u05d5\u05ea \u05e2\u05d3\u05db\u05d5\u05e0\u05d9 \u05de\u05d3\u05d9\u05e0\u05d9\u05d5\u05ea:\",\n  \"\u05e9\u05d9\u05d8\u05ea PPO \u05de\u05e9\u05ea\u05de\u05e9 \u05d1\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05dc\u05de\u05e0\u05d9\u05e2\u05ea \u05e2\u05d3\u05db\u05d5\u05e0\u05d9 \u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05e9\u05e2\u05dc\u05d5\u05dc\u05d9\u05dd \u05dc\u05e2\u05e8\u05e2\u05e8 \u05d0\u05ea \u05d4\u05dc\u05de\u05d9\u05d3\u05d4. \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05d6\u05d4 \u05de\u05d9\u05d5\u05ea\u05e8 \u05e2\u05d1\u05d5\u05e8 RLHF, \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05de\u05e9\u05d8\u05d7 \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d9\u05e6\u05d9\u05d1 \u05d4\u05d5\u05d3\u05d5\u05ea \u05dc warm-started LLM. \u05d4\u05e1\u05e8\u05ea \u05d4\u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05d1-PPO \u05d0\u05d5 \u05d4\u05d9\u05de\u05e0\u05e2\u05d5\u05ea \u05de\u05de\u05e0\u05d5 \u05dc\u05d7\u05dc\u05d5\u05d8\u05d9\u05df \u05e2\u05dd REINFORCE \u05de\u05d5\u05d1\u05d9\u05dc\u05d4 \u05dc\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8, \u05de\u05d4 \u05e9\u05de\u05e6\u05d1\u05d9\u05e2 \u05e2\u05dc \u05db\u05da \u05e9-RLHF \u05d0\u05d9\u05e0\u05d5 \u05d3\u05d5\u05e8\u05e9 \u05e8\u05de\u05d4 \u05db\u05d6\u05d5 \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d1.\",\n  \"4. \u05d0\u05d9\u05d6\u05d5\u05df \u05d1\u05d9\u05df \u05d4\u05e4\u05d7\u05ea\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05d5\u05d4\u05e2\u05dc\u05d0\u05d4 \u05e7\u05dc\u05d4 \u05d1\u05d4\u05d8\u05d9\u05d4:\",\n  \"\u05d0\u05d5\u05de\u05d3\u05df \u05d4\u05d9\u05ea\u05e8\u05d5\u05df \u05d1- PPO \u05de\u05d0\u05d6\u05df \u05d1\u05d9\u05df \u05e9\u05d5\u05e0\u05d5\u05ea \u05d5\u05d4\u05d8\u05d9\u05d4, \u05e0\u05e9\u05dc\u05d8 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d4\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8 \u03bb. \u05e2\u05e8\u05db\u05d9 \u03bb \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 (\u05e7\u05e8\u05d5\u05d1\u05d9\u05dd \u05dc-1) \u05de\u05e4\u05d7\u05d9\u05ea\u05d9\u05dd \u05d4\u05d8\u05d9\u05d4 \u05d0\u05da \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05e9\u05d5\u05e0\u05d5\u05ea. \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05de\u05d9\u05dd \u05e9\u05d1-RLHF, \u05e2\u05e8\u05db\u05d9 \u03bb \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d5\u05d1\u05d9\u05dc\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05e2\u05e7\u05d1\u05d9 \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc, \u05ea\u05d5\u05de\u05db\u05d9\u05dd \u05d1\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05d0\u05d5\u05de\u05d3\u05e0\u05d9\u05dd \u05dc\u05dc\u05d0 \u05d4\u05d8\u05d9\u05d4 \u05db\u05de\u05d5 REINFORCE\",\n  \"\\u000b\u05de\u05d2\u05d1\u05dc\u05d5\u05ea \u05d5\u05db\u05d9\u05d5\u05d5\u05e0\u05d9\u05dd \u05e2\u05ea\u05d9\u05d3\u05d9\u05d9\u05dd\",\n  \"1. \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05d9\u05ea\u05e8 \u05e9\u05dc \u05ea\u05d2\u05de\u05d5\u05dc:\",\n  \"\u05d4\u05de\u05d7\u05e7\u05e8 \u05d0\u05d9\u05e0\u05d5 \u05de\u05ea\u05de\u05d5\u05d3\u05d3 \u05e2\u05dd \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05d9\u05ea\u05e8 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05d4\u05ea\u05d2\u05de\u05d5\u05dc(reward hacking), \u05d1\u05d4 \u05d4\u05de\u05d3\u05d9\u05e0\u05d9\u05d5\u05ea \u05de\u05e0\u05e6\u05dc\u05ea \u05d4\u05d8\u05d9\u05d5\u05ea \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc \u05e2\u05dc \u05d7\u05e9\u05d1\u05d5\u05df \u05d4\u05db\u05dc\u05dc\u05d4. \u05d6\u05d4 \u05e0\u05e9\u05d0\u05e8 \u05d0\u05ea\u05d2\u05e8 \u05e4\u05ea\u05d5\u05d7 \u05e2\u05d1\u05d5\u05e8 RLHF.\",\n  \"2. \u05d4\u05e2\u05e8\u05db\u05d4 \u05d0\u05e0\u05d5\u05e9\u05d9\u05ea:\",\n  \"\u05d1\u05e2\u05d5\u05d3 \u05e9\u05d0\u05d7\u05d5\u05d6\u05d9 \u05e0\u05d9\u05e6\u05d7\u05d5\u05df \u05de\u05d3\u05d5\u05de\u05d9\u05dd \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea GPT-4 \u05de\u05e9\u05de\u05e9\u05d9\u05dd \u05db\u05de\u05d3\u05d3 \u05dc\u05d4\u05e2\u05d3\u05e4\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea, \u05d4\u05e2\u05e8\u05db\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea \u05d9\u05e9\u05d9\u05e8\u05d5\u05ea \u05d4\u05d9\u05d5 \u05de\u05e1\u05e4\u05e7\u05d5\u05ea \u05e8\u05d0\u05d9\u05d5\u05ea \u05d7\u05d6\u05e7\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05dc\u05d0\u05d9\u05db\u05d5\u05ea \u05d4\u05d4\u05ea\u05d0\u05de\u05d4.\",\n  \"3. \\\"\u05e1\u05e7\u05dc\u05d1\u05d9\u05dc\u05d9\u05d5\u05ea\\\":\",\n  \"\u05d4\u05e1\u05e7\u05dc\u05d1\u05d9\u05dc\u05d9\u05d5\u05ea \u05e9\u05dc REINFORCE \u05d5-RLOO \u05dc\u05de\u05d5\u05d3\u05dc\u05d9\u05dd(\u05d4\u05dd \u05d1\u05d3\u05e7\u05d5 \u05e8\u05e7 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05dc 7B) \u05d5\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05e6\u05e8\u05d9\u05db\u05d4 \u05de\u05d7\u05e7\u05e8 \u05e0\u05d5\u05e1\u05e3.\",\n  \"\u05de\u05e1\u05e7\u05e0\u05d4\",\n  \"\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05d2 \u05d8\u05d9\u05e2\u05d5\u05df \u05de\u05e9\u05db\u05e0\u05e2 \u05dc\u05d1\u05d7\u05d9\u05e0\u05d4 \u05de\u05d7\u05d5\u05d3\u05e9\u05ea \u05e9\u05dc \u05e9\u05d9\u05d8\u05d5\u05ea \u05d1\u05e1\u05d2\u05e0\u05d5\u05df REINFORCE \u05d1-RLHF, \u05de\u05d0\u05ea\u05d2\u05e8 \u05d0\u05ea \u05d4\u05d3\u05d5\u05de\u05d9\u05e0\u05e0\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc   PPO \u05d5\u05d3\u05d5\u05de\u05d9\u05d4. \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e0\u05d9\u05e6\u05d5\u05dc \u05d4\u05de\u05d0\u05e4\u05d9\u05d9\u05e0\u05d9\u05dd \u05d4\u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d9\u05dd \u05e9\u05dc RLHF - \u05db\u05de\u05d5 warm started LLM \u05d5\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05d1\u05e8\u05de\u05ea \u05d4\u05e1\u05d3\u05e8\u05d4 - \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05de\u05d9\u05dd \u05e9\u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05db\u05de\u05d5 REINFORCE \u05d5-RLOO \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05dc\u05e2\u05dc\u05d5\u05ea \u05e2\u05dc \u05d7\u05dc\u05d5\u05e4\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05db\u05de\u05d5 PPO \u05d5-RAFT \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc, \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05d5\u05e2\u05de\u05d9\u05d3\u05d5\u05ea.\",\n  \"https://arxiv.org/abs/2402.14740\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst originalCount = paragraphs.items.length; // 6 in before.docx\n\n// 1) Update the text of the paragraphs that already exist (indices 0..originalCount-1)\nfor (let i = 0; i < originalCount; i++) {\n  paragraphs.items[i].insertText(targetTexts[i], Word.InsertLocation.replace);\n}\nawait context.sync();\n\n// 2) Append the remaining new paragraphs after the last existing paragraph,\n//    preserving order.\nlet anchor = paragraphs.items[originalCount - 1];\nfor (let i = originalCount; i < targetTexts.length; i++) {\n  anchor = anchor.insertParagraph(targetTexts[i], Word.InsertLocation.after);\n}\nawait context.sync();\n\nreturn \"ok: \" + targetTexts.length + \" paragraphs\";\n", "ps1": "# Target paragraph texts after the edit (28 paragraphs total).\n# [char]11 (vertical tab) represents an in-paragraph line break (<w:br/>),\n# matching how Word exposes a manual line break through Range.Text.\n$targetTexts = @(\n  ('\u05d4\u05de\u05d0\u05de\u05e8 \u05d4\u05d9\u05d5\u05de\u05d9 \u05e9\u05dc \u05de\u05d9\u05d9\u05e7 - 24.01.25'+[char]11+'Back to Basics: Revisiting REINFORCE Style Optimization for Learning from Human Feedback in LLMs'),\n  ('\u05ea\u05de\u05e6\u05d9\u05ea \u05d4\u05de\u05d0\u05de\u05e8:'),\n  ('\u05d4\u05de\u05d0\u05de\u05e8 \u05d1\u05d5\u05d7\u05df \u05de\u05d7\u05d3\u05e9 \u05d0\u05ea \u05d4\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05dc\u05de\u05d9\u05d3\u05d4 \u05de\u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd \u05de\u05e4\u05d9\u05d3\u05d1\u05e7 \u05d0\u05e0\u05d5\u05e9\u05d9 (RLHF) \u05d1\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05ea LLMs. \u05d4\u05d5\u05d0 \u05de\u05d0\u05ea\u05d2\u05e8 \u05d0\u05ea \u05d4\u05d3\u05d5\u05de\u05d9\u05e0\u05e0\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc PPO (Proximal Policy Optimization) \u05db\u05e9\u05d9\u05d8\u05ea \u05dc\u05de\u05d9\u05d3\u05ea \u05d4\u05d7\u05d9\u05d6\u05d5\u05e7\u05d9\u05dd \u05d4\u05e1\u05d8\u05e0\u05d3\u05e8\u05d8\u05d9\u05ea \u05d1\u05d4\u05e7\u05e9\u05e8 \u05d6\u05d4, \u05ea\u05d5\u05da \u05d4\u05d3\u05d2\u05e9\u05ea \u05d7\u05d5\u05e1\u05e8 \u05d4\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d4\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea \u05d5\u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d4\u05de\u05d9\u05d5\u05ea\u05e8\u05ea \u05e9\u05dc\u05d5. \u05d1\u05de\u05e7\u05d5\u05dd \u05d6\u05d0\u05ea, \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e6\u05d9\u05e2\u05d9\u05dd \u05dc\u05d7\u05d6\u05d5\u05e8 \u05dc\u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d1\u05e1\u05d2\u05e0\u05d5\u05df REINFORCE, \u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05ea Vanilla Policy Gradient (REINFORCE) \u05d5\u05d4\u05e8\u05d7\u05d1\u05ea\u05d5 \u05d4\u05e8\u05d1-\u05d3\u05d2\u05d9\u05de\u05ea\u05d9\u05ea, REINFORCE Leave-One-Out (RLOO). \u05e9\u05d9\u05d8\u05d5\u05ea \u05d0\u05dc\u05d5 \u05de\u05d5\u05db\u05d9\u05d7\u05d5\u05ea \u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de-PPO \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05e2\u05dc\u05d5\u05ea \u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05ea, \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05d1\u05de\u05e1\u05e4\u05e8 \u05de\u05e2\u05e8\u05db\u05d9 \u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05d5\u05d0\u05e8\u05db\u05d9\u05d8\u05e7\u05d8\u05d5\u05e8\u05d5\u05ea LLM. \u05d4\u05de\u05de\u05e6\u05d0\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05e9\u05d9\u05dd \u05e9\u05e0\u05d9\u05ea\u05df \u05dc\u05d4\u05e9\u05d9\u05d2 \u05d4\u05ea\u05d0\u05de\u05d4 \u05e9\u05dc LLMs \u05dc\u05d4\u05e2\u05d3\u05e4\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea \u05e2\u05dd \u05d0\u05e1\u05d8\u05e8\u05d8\u05d2\u05d9\u05d5\u05ea \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05d4\u05de\u05d5\u05ea\u05d0\u05de\u05d5\u05ea \u05dc\u05d9\u05d9\u05d7\u05d5\u05d3\u05d9\u05d5\u05ea \u05e9\u05dc RLHF.'),\n  ('\u05d4\u05e8\u05d7\u05d1\u05d4 \u05e2\u05dc \u05e0\u05e7\u05d5\u05d3\u05d5\u05ea \u05e2\u05d9\u05e7\u05e8\u05d9\u05d5\u05ea:'),\n  ('1. \u05e4\u05d9\u05e9\u05d5\u05d8 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9:'),\n  ('\u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05e8\u05d1\u05d9\u05dd \u05de\u05d4\u05e8\u05db\u05d9\u05d1\u05d9\u05dd \u05e9\u05dc PPO (\u05dc\u05de\u05e9\u05dc, \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2, \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05e8\u05da (value), \u05d5\u05de\u05d9\u05d3\u05d5\u05dc  \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd) \u05d0\u05d9\u05e0\u05dd \u05d4\u05db\u05e8\u05d7\u05d9\u05d9\u05dd \u05dc-RLHF, \u05d1\u05d4\u05d9\u05e0\u05ea\u05df \u05d4\u05d0\u05ea\u05d7\u05d5\u05dc \u05d8\u05d5\u05d1 \u05e9\u05dc LLM (\u05dc\u05d0\u05d7\u05e8 SFT \u05dc\u05de\u05e9\u05dc). \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d9\u05d3\u05d5\u05dc \u05e1\u05d3\u05e8\u05d5\u05ea \u05e9\u05dc\u05de\u05d5\u05e5 \u05db\u05e4\u05e2\u05d5\u05dc\u05d5\u05ea \u05d1\u05d5\u05d3\u05d3\u05d5\u05ea, REINFORCE \u05e0\u05de\u05e0\u05e2 \u05de\u05d4\u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05e9\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05e8\u05da-\u05de\u05e6\u05d1(V \u05d5-Q) \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05df, \u05d5\u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05d1\u05e2\u05d9\u05d4 \u05dc\u05d3\u05d5\u05de\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05dc\u05d1\u05e0\u05d3\u05d9\u05d8 \u05d4\u05e7\u05e9\u05e8\u05d9.'),\n  ('2. \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05de\u05e2\u05e9\u05d9\u05ea:'),\n  ('\u05e9\u05d9\u05d8\u05ea RLOO \u05de\u05e9\u05ea\u05de\u05e9\u05ea \u05d1\u05db\u05dc \u05d4\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05e9\u05e0\u05d5\u05e6\u05e8\u05d5 \u05dc\u05d1\u05e0\u05d9\u05d9\u05ea \u05d1\u05e1\u05d9\u05e1 \u05d4\u05e9\u05d5\u05d5\u05d0\u05d4, \u05de\u05e9\u05d9\u05d2 \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05d2\u05d1\u05d5\u05d4\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de-RAFT, \u05e9\u05dc\u05d5\u05e7\u05d7 \u05d0\u05ea \u05e8\u05e7 \u05d0\u05ea \u05d4\u05d3\u05d2\u05d9\u05de\u05d5\u05ea \u05d1\u05e2\u05dc\u05d5\u05ea \u05e6\u05d9\u05d5\u05e0\u05d9\u05dd \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd (\u05e1\u05d5\u05d2 \u05e9\u05dc rejection sampling). \u05d6\u05d4 \u05de\u05d5\u05d1\u05d9\u05dc \u05dc\u05d7\u05d9\u05e1\u05db\u05d5\u05df \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9 \u05d1\u05d7\u05d9\u05e9\u05d5\u05d1\u05d9\u05dd \u05d5\u05e0\u05d9\u05e6\u05d5\u05dc \u05d8\u05d5\u05d1 \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05d4\u05e0\u05ea\u05d5\u05e0\u05d9\u05dd \u05d4\u05d6\u05de\u05d9\u05e0\u05d9\u05dd. \u05d4\u05d2\u05d9\u05e9\u05d4 \u05de\u05e4\u05e9\u05d8\u05ea \u05d0\u05ea \u05ea\u05d4\u05dc\u05d9\u05db\u05d9 \u05d4-RLHF \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05e4\u05d7\u05ea\u05ea \u05d4\u05ea\u05dc\u05d5\u05ea \u05d1\u05d4\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05e8\u05d2\u05d9\u05e9\u05d9\u05dd \u05db\u05de\u05d5 \u05d9\u05d7\u05e1\u05d9 \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05d5\u05d4\u05e4\u05e8\u05de\u05d8\u05e8\u05d9\u05dd \u05d1\u05e9\u05e2\u05e8\u05d5\u05da \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d9\u05ea\u05e8\u05d5\u05df (\u05db\u05de\u05d5 \u05d1-GAE).'),\n  ('4.\u05e8\u05d5\u05d1\u05e1\u05d8\u05d9\u05d5\u05ea:'+[char]11+'\u05e9\u05d9\u05d8\u05ea RLOO  \u05de\u05d3\u05d2\u05d9\u05de\u05d4 \u05e8\u05d5\u05d1\u05e1\u05d8\u05d9\u05d5\u05ea \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05e8\u05d5\u05e2\u05e9\u05d9\u05dd \u05d5\u05e2\u05d5\u05e0\u05e9\u05d9 KL \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8, \u05e2\u05d5\u05dc\u05d4 \u05e2\u05dc \u05e9\u05d9\u05d8\u05d5\u05ea \u05db\u05de\u05d5 RAFT \u05e9\u05e8\u05d2\u05d9\u05e9\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05dc\u05d3\u05d9\u05d5\u05e7\u05dd.'),\n  ('\u05ea\u05d5\u05d1\u05e0\u05d5\u05ea \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05d5\u05ea:'),\n  ('1. \u05d0\u05d9\u05d6\u05d5\u05df \u05e9\u05d5\u05e0\u05d5\u05ea-\u05d4\u05d8\u05d9\u05d4(bias-variance tradeoff) \u05d5\u05e9\u05e2\u05e8\u05d5\u05da \u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05dc\u05dc\u05d0 \u05d4\u05d8\u05d9\u05d4:'),\n  ('\u05e9\u05d9\u05d8\u05ea PPO \u05de\u05e1\u05ea\u05de\u05db\u05ea \u05e2\u05dc \u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d5\u05ea \u05e2\u05e8\u05da-\u05de\u05e6\u05d1 \u05d5\u05e9\u05e2\u05e8\u05d5\u05da \u05d9\u05ea\u05e8\u05d5\u05df \u05de\u05d5\u05db\u05dc\u05dc (Generalized Advantage Estimation) \u05dc\u05d4\u05e4\u05d7\u05ea\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05e9\u05dc \u05e9\u05e2\u05e8\u05d5\u05da \u05d2\u05e8\u05d3\u05d9\u05d0\u05e0\u05d8 \u05d1\u05de\u05d7\u05d9\u05e8 \u05e9\u05dc \u05d4\u05e2\u05dc\u05d0\u05ea \u05d4\u05d8\u05d9\u05d4. \u05d4\u05de\u05d0\u05de\u05e8 \u05d8\u05d5\u05e2\u05df \u05e9\u05d1-RLHF \u05e2\u05d1\u05d5\u05e8 LLMs \u05d3\u05d9 \u05de\u05d0\u05d5\u05de\u05df (warm start) \u05d4\u05d5\u05e4\u05da \u05d0\u05ea \u05d4\u05e4\u05d7\u05ea\u05ea \u05d4\u05e9\u05d5\u05e0\u05d5\u05ea \u05dc\u05e4\u05d7\u05d5\u05ea \u05e7\u05e8\u05d9\u05d8\u05d9\u05ea. \u05d6\u05d4 \u05de\u05d0\u05e4\u05e9\u05e8 \u05dc\u05e9\u05d9\u05d8\u05d5\u05ea \u05dc\u05dc\u05d0 \u05d4\u05d8\u05d9\u05d4 \u05db\u05de\u05d5 REINFORCE \u05dc\u05ea\u05e4\u05e7\u05d3 \u05d4\u05d9\u05d8\u05d1 \u05d1\u05dc\u05d9 \u05dc\u05d4\u05db\u05e0\u05d9\u05e1 \u05d4\u05d8\u05d9\u05d4 \u05de\u05e9\u05de\u05e2\u05d5\u05ea\u05d9\u05ea. \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05ea, \u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05d3\u05d2\u05d9\u05dd \u05e9-REINFORCE \u05de\u05e9\u05d9\u05d2 \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc \u05d8\u05d5\u05d1\u05d4 \u05d9\u05d5\u05ea\u05e8 \u05de-PPO, \u05d0\u05e4\u05d9\u05dc\u05d5 \u05ea\u05d7\u05ea \u05ea\u05e0\u05d0\u05d9\u05dd \u05e9\u05dc \u05e9\u05d5\u05e0\u05d5\u05ea \u05d2\u05d1\u05d5\u05d4\u05d4 \u05ea\u05d9\u05d0\u05d5\u05e8\u05d8\u05d9\u05ea.'),\n  ('2. \u05de\u05d9\u05d3\u05d5\u05dc \u05de\u05e1\u05dc\u05d5\u05dc \u05de\u05dc\u05d0(\u05ea\u05e9\u05d5\u05d1\u05d4 \u05e9\u05dc\u05de\u05d4) \u05dc\u05e2\u05d5\u05de\u05ea \u05de\u05d9\u05d3\u05d5\u05dc \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd:'),\n  ('\u05e9\u05d9\u05d8\u05ea PPO \u05de\u05de\u05d3\u05dc\u05ea \u05db\u05dc \u05d8\u05d5\u05e7\u05df \u05db\u05e4\u05e2\u05d5\u05dc\u05d4, \u05d9\u05d5\u05e6\u05e8 \u05ea\u05d4\u05dc\u05d9\u05da \u05d4\u05d7\u05dc\u05d8\u05d4 \u05de\u05e8\u05e7\u05d5\u05d1\u05d9 (MDP) \u05d1\u05d5 \u05e8\u05e6\u05e4\u05d9 \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d7\u05dc\u05e7\u05d9\u05d9\u05dd \u05d4\u05dd \u05de\u05e6\u05d1\u05d9\u05dd. \u05e2\u05dd \u05d6\u05d0\u05ea, RLHF \u05de\u05d9\u05d9\u05d7\u05e1 \u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05e8\u05e7 \u05dc\u05ea\u05e9\u05d5\u05d1\u05d5\u05ea \u05e9\u05dc\u05de\u05d5\u05ea, \u05de\u05d4 \u05e9\u05d4\u05d5\u05e4\u05da \u05de\u05e6\u05d1\u05d9 \u05d1\u05d9\u05e0\u05d9\u05d9\u05dd \u05dc\u05dc\u05d0 \u05e8\u05dc\u05d5\u05d5\u05e0\u05d8\u05d9\u05d9\u05dd. \u05e2\u05dc \u05d9\u05d3\u05d9 \u05de\u05d9\u05d3\u05d5\u05dc \u05d4\u05ea\u05e9\u05d5\u05d1\u05d4 \u05db\u05e4\u05e2\u05d5\u05dc\u05d4 \u05d9\u05d7\u05d9\u05d3\u05d4, REINFORCE \u05de\u05e4\u05e9\u05d8 \u05d0\u05ea \u05d4\u05d1\u05e2\u05d9\u05d4 \u05dc\u05de\u05d1\u05e0\u05d4 Contextual Bandit \u05d4\u05de\u05ea\u05d9\u05d9\u05e9\u05e8 \u05d9\u05e9\u05d9\u05e8\u05d5\u05ea \u05e2\u05dd \u05de\u05d1\u05e0\u05d4 \u05d4\u05ea\u05d2\u05de\u05d5\u05dc. \u05ea\u05d5\u05e6\u05d0\u05d5\u05ea \u05d0\u05de\u05e4\u05d9\u05e8\u05d9\u05d5\u05ea \u05de\u05d0\u05e9\u05e8\u05d5\u05ea \u05db\u05d9 \u05d2\u05d9\u05e9\u05d4 \u05d6\u05d5 \u05e2\u05d5\u05dc\u05d4 \u05e2\u05dc \u05de\u05d9\u05d3\u05d5\u05dc \u05d1\u05e8\u05de\u05ea \u05d8\u05d5\u05e7\u05e0\u05d9\u05dd \u05d4\u05df \u05d1\u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d5\u05d4\u05df \u05d1\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd.'),\n  ('3. \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05d5\u05d9\u05e6\u05d9\u05d1\u05d5\u05ea \u05e2\u05d3\u05db\u05d5\u05e0\u05d9 \u05de\u05d3\u05d9\u05e0\u05d9\u05d5\u05ea:'),\n  ('\u05e9\u05d9\u05d8\u05ea PPO \u05de\u05e9\u05ea\u05de\u05e9 \u05d1\u05de\u05e0\u05d2\u05e0\u05d5\u05df \u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05dc\u05de\u05e0\u05d9\u05e2\u05ea \u05e2\u05d3\u05db\u05d5\u05e0\u05d9 \u05e4\u05d5\u05dc\u05d9\u05e1\u05d9 \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05e9\u05e2\u05dc\u05d5\u05dc\u05d9\u05dd \u05dc\u05e2\u05e8\u05e2\u05e8 \u05d0\u05ea \u05d4\u05dc\u05de\u05d9\u05d3\u05d4. \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05e8\u05d0\u05d9\u05dd \u05e9\u05d6\u05d4 \u05de\u05d9\u05d5\u05ea\u05e8 \u05e2\u05d1\u05d5\u05e8 RLHF, \u05de\u05db\u05d9\u05d5\u05d5\u05df \u05e9\u05de\u05e9\u05d8\u05d7 \u05d4\u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d4 \u05d9\u05e6\u05d9\u05d1 \u05d4\u05d5\u05d3\u05d5\u05ea \u05dc warm-started LLM. \u05d4\u05e1\u05e8\u05ea \u05d4\u05e7\u05dc\u05d9\u05e4\u05d9\u05e0\u05d2 \u05d1-PPO \u05d0\u05d5 \u05d4\u05d9\u05de\u05e0\u05e2\u05d5\u05ea \u05de\u05de\u05e0\u05d5 \u05dc\u05d7\u05dc\u05d5\u05d8\u05d9\u05df \u05e2\u05dd REINFORCE \u05de\u05d5\u05d1\u05d9\u05dc\u05d4 \u05dc\u05d1\u05d9\u05e6\u05d5\u05e2\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8, \u05de\u05d4 \u05e9\u05de\u05e6\u05d1\u05d9\u05e2 \u05e2\u05dc \u05db\u05da \u05e9-RLHF \u05d0\u05d9\u05e0\u05d5 \u05d3\u05d5\u05e8\u05e9 \u05e8\u05de\u05d4 \u05db\u05d6\u05d5 \u05e9\u05dc \u05d9\u05d9\u05e6\u05d5\u05d1.'),\n  ('4. \u05d0\u05d9\u05d6\u05d5\u05df \u05d1\u05d9\u05df \u05d4\u05e4\u05d7\u05ea\u05ea \u05e9\u05d5\u05e0\u05d5\u05ea \u05d5\u05d4\u05e2\u05dc\u05d0\u05d4 \u05e7\u05dc\u05d4 \u05d1\u05d4\u05d8\u05d9\u05d4:'),\n  ('\u05d0\u05d5\u05de\u05d3\u05df \u05d4\u05d9\u05ea\u05e8\u05d5\u05df \u05d1- PPO \u05de\u05d0\u05d6\u05df \u05d1\u05d9\u05df \u05e9\u05d5\u05e0\u05d5\u05ea \u05d5\u05d4\u05d8\u05d9\u05d4, \u05e0\u05e9\u05dc\u05d8 \u05e2\u05dc \u05d9\u05d3\u05d9 \u05d4\u05d4\u05d9\u05e4\u05e8-\u05e4\u05e8\u05de\u05d8\u05e8 \u03bb. \u05e2\u05e8\u05db\u05d9 \u03bb \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 (\u05e7\u05e8\u05d5\u05d1\u05d9\u05dd \u05dc-1) \u05de\u05e4\u05d7\u05d9\u05ea\u05d9\u05dd \u05d4\u05d8\u05d9\u05d4 \u05d0\u05da \u05de\u05d2\u05d3\u05d9\u05dc\u05d9\u05dd \u05e9\u05d5\u05e0\u05d5\u05ea. \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05de\u05d9\u05dd \u05e9\u05d1-RLHF, \u05e2\u05e8\u05db\u05d9 \u03bb \u05d2\u05d1\u05d5\u05d4\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05d5\u05d1\u05d9\u05dc\u05d9\u05dd \u05d1\u05d0\u05d5\u05e4\u05df \u05e2\u05e7\u05d1\u05d9 \u05dc\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05d8\u05d5\u05d1\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05e9\u05dc \u05d4\u05de\u05d5\u05d3\u05dc, \u05ea\u05d5\u05de\u05db\u05d9\u05dd \u05d1\u05e9\u05d9\u05de\u05d5\u05e9 \u05d1\u05d0\u05d5\u05de\u05d3\u05e0\u05d9\u05dd \u05dc\u05dc\u05d0 \u05d4\u05d8\u05d9\u05d4 \u05db\u05de\u05d5 REINFORCE'),\n  ([char]11+'\u05de\u05d2\u05d1\u05dc\u05d5\u05ea \u05d5\u05db\u05d9\u05d5\u05d5\u05e0\u05d9\u05dd \u05e2\u05ea\u05d9\u05d3\u05d9\u05d9\u05dd'),\n  ('1. \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05d9\u05ea\u05e8 \u05e9\u05dc \u05ea\u05d2\u05de\u05d5\u05dc:'),\n  ('\u05d4\u05de\u05d7\u05e7\u05e8 \u05d0\u05d9\u05e0\u05d5 \u05de\u05ea\u05de\u05d5\u05d3\u05d3 \u05e2\u05dd \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05d9\u05ea\u05e8 \u05e9\u05dc \u05de\u05d5\u05d3\u05dc \u05d4\u05ea\u05d2\u05de\u05d5\u05dc(reward hacking), \u05d1\u05d4 \u05d4\u05de\u05d3\u05d9\u05e0\u05d9\u05d5\u05ea \u05de\u05e0\u05e6\u05dc\u05ea \u05d4\u05d8\u05d9\u05d5\u05ea \u05d1\u05e4\u05d5\u05e0\u05e7\u05e6\u05d9\u05d9\u05ea \u05d4\u05ea\u05d2\u05de\u05d5\u05dc \u05e2\u05dc \u05d7\u05e9\u05d1\u05d5\u05df \u05d4\u05db\u05dc\u05dc\u05d4. \u05d6\u05d4 \u05e0\u05e9\u05d0\u05e8 \u05d0\u05ea\u05d2\u05e8 \u05e4\u05ea\u05d5\u05d7 \u05e2\u05d1\u05d5\u05e8 RLHF.'),\n  ('2. \u05d4\u05e2\u05e8\u05db\u05d4 \u05d0\u05e0\u05d5\u05e9\u05d9\u05ea:'),\n  ('\u05d1\u05e2\u05d5\u05d3 \u05e9\u05d0\u05d7\u05d5\u05d6\u05d9 \u05e0\u05d9\u05e6\u05d7\u05d5\u05df \u05de\u05d3\u05d5\u05de\u05d9\u05dd \u05d1\u05d0\u05de\u05e6\u05e2\u05d5\u05ea GPT-4 \u05de\u05e9\u05de\u05e9\u05d9\u05dd \u05db\u05de\u05d3\u05d3 \u05dc\u05d4\u05e2\u05d3\u05e4\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea, \u05d4\u05e2\u05e8\u05db\u05d5\u05ea \u05d0\u05e0\u05d5\u05e9\u05d9\u05d5\u05ea \u05d9\u05e9\u05d9\u05e8\u05d5\u05ea \u05d4\u05d9\u05d5 \u05de\u05e1\u05e4\u05e7\u05d5\u05ea \u05e8\u05d0\u05d9\u05d5\u05ea \u05d7\u05d6\u05e7\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05dc\u05d0\u05d9\u05db\u05d5\u05ea \u05d4\u05d4\u05ea\u05d0\u05de\u05d4.'),\n  ('3. \"\u05e1\u05e7\u05dc\u05d1\u05d9\u05dc\u05d9\u05d5\u05ea\":'),\n  ('\u05d4\u05e1\u05e7\u05dc\u05d1\u05d9\u05dc\u05d9\u05d5\u05ea \u05e9\u05dc REINFORCE \u05d5-RLOO \u05dc\u05de\u05d5\u05d3\u05dc\u05d9\u05dd(\u05d4\u05dd \u05d1\u05d3\u05e7\u05d5 \u05e8\u05e7 \u05de\u05d5\u05d3\u05dc\u05d9\u05dd \u05e9\u05dc 7B) \u05d5\u05d3\u05d0\u05d8\u05d4\u05e1\u05d8\u05d9\u05dd \u05d2\u05d3\u05d5\u05dc\u05d9\u05dd \u05d9\u05d5\u05ea\u05e8 \u05de\u05e6\u05e8\u05d9\u05db\u05d4 \u05de\u05d7\u05e7\u05e8 \u05e0\u05d5\u05e1\u05e3.'),\n  ('\u05de\u05e1\u05e7\u05e0\u05d4'),\n  ('\u05d4\u05de\u05d0\u05de\u05e8 \u05de\u05e6\u05d9\u05d2 \u05d8\u05d9\u05e2\u05d5\u05df \u05de\u05e9\u05db\u05e0\u05e2 \u05dc\u05d1\u05d7\u05d9\u05e0\u05d4 \u05de\u05d7\u05d5\u05d3\u05e9\u05ea \u05e9\u05dc \u05e9\u05d9\u05d8\u05d5\u05ea \u05d1\u05e1\u05d2\u05e0\u05d5\u05df REINFORCE \u05d1-RLHF, \u05de\u05d0\u05ea\u05d2\u05e8 \u05d0\u05ea \u05d4\u05d3\u05d5\u05de\u05d9\u05e0\u05e0\u05d8\u05d9\u05d5\u05ea \u05e9\u05dc   PPO \u05d5\u05d3\u05d5\u05de\u05d9\u05d4. \u05e2\u05dc \u05d9\u05d3\u05d9 \u05e0\u05d9\u05e6\u05d5\u05dc \u05d4\u05de\u05d0\u05e4\u05d9\u05d9\u05e0\u05d9\u05dd \u05d4\u05e1\u05e4\u05e6\u05d9\u05e4\u05d9\u05d9\u05dd \u05e9\u05dc RLHF - \u05db\u05de\u05d5 warm started LLM \u05d5\u05ea\u05d2\u05de\u05d5\u05dc\u05d9\u05dd \u05d1\u05e8\u05de\u05ea \u05d4\u05e1\u05d3\u05e8\u05d4 - \u05d4\u05d7\u05d5\u05e7\u05e8\u05d9\u05dd \u05de\u05d3\u05d2\u05d9\u05de\u05d9\u05dd \u05e9\u05e9\u05d9\u05d8\u05d5\u05ea \u05e4\u05e9\u05d5\u05d8\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05db\u05de\u05d5 REINFORCE \u05d5-RLOO \u05d9\u05db\u05d5\u05dc\u05d5\u05ea \u05dc\u05e2\u05dc\u05d5\u05ea \u05e2\u05dc \u05d7\u05dc\u05d5\u05e4\u05d5\u05ea \u05de\u05d5\u05e8\u05db\u05d1\u05d5\u05ea \u05d9\u05d5\u05ea\u05e8 \u05db\u05de\u05d5 PPO \u05d5-RAFT \u05de\u05d1\u05d7\u05d9\u05e0\u05ea \u05d0\u05d5\u05e4\u05d8\u05d9\u05de\u05d9\u05d6\u05e6\u05d9\u05d9\u05ea \u05ea\u05d2\u05de\u05d5\u05dc, \u05d9\u05e2\u05d9\u05dc\u05d5\u05ea \u05d3\u05d2\u05d9\u05de\u05d4 \u05d5\u05e2\u05de\u05d9\u05d3\u05d5\u05ea.'),\n  ('https://arxiv.org/abs/2402.14740')\n)\n\n$d = $word.ActiveDocument\n\n$originalCount = $d.Paragraphs.Count   # 6 in before.docx\n\n# 1) Overwrite the text of the paragraphs that already exist. Assigning to\n#    Range.Text replaces the paragraph's run content in place (keeps the\n#    paragraph mark / pPr, including the \"Normal\" style, untouched).\nfor ($i = 1; $i -le $originalCount; $i++) {\n    $d.Paragraphs($i).Range.Text = $targetTexts[$i - 1]\n}\n\n# 2) Append the remaining new paragraphs, one at a time, after the last\n#    existing paragraph, preserving order.\nfor ($i = $originalCount; $i -lt $targetTexts.Count; $i++) {\n    $d.Paragraphs($d.Paragraphs.Count).Range.InsertParagraphAfter()\n    $d.Paragraphs($d.Paragraphs.Count).Range.Text = $targetTexts[$i]\n}\n\nWrite-Output (\"done: \" + $d.Paragraphs.Count + \" paragraphs\")\n"}
